# Applies the diff to Payments sheet:
#  - Adds new column G "IsDeleted" (header styled like other headers)
#  - Updates values in rows 2-4 for columns B, C, D, E, F
#  - Adds 0 values in new column G for rows 2-4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell G1, copying the style from the existing header F1
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "IsDeleted"

# Row 2
$ws.Range("B2").Value = 16
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 9
$ws.Range("E2").Value = 222
$ws.Range("F2").Value = "2025-03-23 21:55:28"
$ws.Range("G2").Value = 0

# Row 3
$ws.Range("B3").Value = 20
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 9
$ws.Range("E3").Value = 920
$ws.Range("F3").Value = "2025-03-23 22:07:30"
$ws.Range("G3").Value = 0

# Row 4
$ws.Range("B4").Value = 16
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 9
$ws.Range("E4").Value = 718
$ws.Range("F4").Value = "2025-03-23 22:14:37"
$ws.Range("G4").Value = 0

Write-Host "Applied Payments sheet updates. UsedRange:" $ws.UsedRange.Address()
